# Fruta / hortaliza, semanal
# Insert two new weekly-price rows (2021-09-09, "amarillo" quality) above the
# existing row 234, shifting the old rows 234-239 down to 236-241.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 234 (pushes existing 234-239 down to 236-241).
$ws.Rows.Item(234).Insert()
$ws.Rows.Item(234).Insert()

# --- Row 234 (new): Limón, 1a amarillo ---
$ws.Cells.Item(234, 1).Value = 11
$ws.Cells.Item(234, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(234, 3).Value = "Bíobío"
$ws.Cells.Item(234, 4).Value = 44448
$ws.Cells.Item(234, 5).Value = 8
$ws.Cells.Item(234, 6).Value = "Fruta"
$ws.Cells.Item(234, 7).Value = 100102
$ws.Cells.Item(234, 8).Value = "Cítricos"
$ws.Cells.Item(234, 9).Value = 100102003
$ws.Cells.Item(234, 10).Value = "Limón"
$ws.Cells.Item(234, 11).Value = "Sin especificar"
$ws.Cells.Item(234, 12).Value = "1a amarillo"
$ws.Cells.Item(234, 13).Value = 300
$ws.Cells.Item(234, 14).Value = 5000
$ws.Cells.Item(234, 15).Value = 5000
$ws.Cells.Item(234, 16).Value = 5000
$ws.Cells.Item(234, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(234, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(234, 19).Value = 312
$ws.Cells.Item(234, 20).Value = 16

# --- Row 235 (new): Limón, 2a amarillo ---
$ws.Cells.Item(235, 1).Value = 11
$ws.Cells.Item(235, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(235, 3).Value = "Bíobío"
$ws.Cells.Item(235, 4).Value = 44448
$ws.Cells.Item(235, 5).Value = 8
$ws.Cells.Item(235, 6).Value = "Fruta"
$ws.Cells.Item(235, 7).Value = 100102
$ws.Cells.Item(235, 8).Value = "Cítricos"
$ws.Cells.Item(235, 9).Value = 100102003
$ws.Cells.Item(235, 10).Value = "Limón"
$ws.Cells.Item(235, 11).Value = "Sin especificar"
$ws.Cells.Item(235, 12).Value = "2a amarillo"
$ws.Cells.Item(235, 13).Value = 300
$ws.Cells.Item(235, 14).Value = 4000
$ws.Cells.Item(235, 15).Value = 4000
$ws.Cells.Item(235, 16).Value = 4000
$ws.Cells.Item(235, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(235, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(235, 19).Value = 250
$ws.Cells.Item(235, 20).Value = 16

Write-Output "done"
